$d = $word.ActiveDocument

# The two empty paragraphs immediately following the "{{cover_sheet_top_message}}"
# and "{{ cover_sheet_share }}" headings currently have no paragraph properties
# (no <w:pPr>), which defaults to full "space after" below them. Reduce the
# space under those headings by giving each paragraph explicit spacing of 0
# points after, matching the neighboring paragraphs' formatting.
$targets = @(6, 14)

foreach ($idx in $targets) {
    $p = $d.Paragraphs($idx)
    if ($p.Range.Text.Trim("`r", "`a") -ne "") {
        throw "Paragraph $idx is not empty as expected"
    }
    $p.SpaceAfter = 0
}
